# Commit: "added SCT rank and -1 to adwDF"
#
# Target change (per the OOXML diff):
#   - ALPHA sheet (sheet3): insert 9 new personnel rows with rank "SCT"
#     right after the "NIL / OC" row and before the "CPT" rows.
#   - ALPHA sheet: rename displayName/sheetName "RICHMOND RAY" -> "RICHMOND"
#     and "LEON LAI" -> "LEON" (both in the 2LT block).
#   - No other sheets have content changes (BRAVO/OTHERS only see shared
#     string index shifts caused by the new strings, not real edits).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALPHA")

# Insert 9 blank rows right before the current row 3 (CPT / MARC),
# i.e. immediately after the NIL / OC row.
$ws.Range("A3:E11").EntireRow.Insert()

# New SCT personnel to add (rank, displayName, sheetName, commSec, nor)
$newPersonnel = @(
    @("SCT", "BRAYDEN",    "BRAYDEN",    "NIL", "REGULAR"),
    @("SCT", "NOAH LAM",   "NOAH LAM",   "NIL", "NSF"),
    @("SCT", "MARCUS",     "MARCUS",     "NIL", "NSF"),
    @("SCT", "MENG LONG",  "MENG LONG",  "NIL", "NSF"),
    @("SCT", "KAI",        "KAI",        "NIL", "NSF"),
    @("SCT", "CHARLES",    "CHARLES",    "NIL", "NSF"),
    @("SCT", "DARSHAN",    "DARSHAN",    "NIL", "NSF"),
    @("SCT", "ZHONG PING", "ZHONG PING", "NIL", "NSF"),
    @("SCT", "DERRILL",    "DERRILL",    "NIL", "NSF")
)

for ($i = 0; $i -lt $newPersonnel.Count; $i++) {
    $row = 3 + $i
    $rec = $newPersonnel[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
}

# Rename the two existing 2LT personnel (now shifted down 9 rows, to 18/19)
$ws.Range("B18").Value = "RICHMOND"
$ws.Range("C18").Value = "RICHMOND"
$ws.Range("B19").Value = "LEON"
$ws.Range("C19").Value = "LEON"

Write-Output "Added $($newPersonnel.Count) SCT personnel to ALPHA; renamed RICHMOND RAY -> RICHMOND, LEON LAI -> LEON"
